$wb = $excel.ActiveWorkbook

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3409445.5
$ws.Range("I86").Value = 43993.332
$ws.Range("J86").Value = 5092171.5
$ws.Range("K86").Value = 43993.332
$ws.Range("L86").Value = 5092171.5
$ws.Range("M86").Value = -42870.332
$ws.Range("N86").Value = -5094417.5

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3409445.5
$ws.Range("I89").Value = 43993.332
$ws.Range("J89").Value = 5092171.5
$ws.Range("K89").Value = 219966.66
$ws.Range("L89").Value = 25460857.5
$ws.Range("M89").Value = -214350.66
$ws.Range("N89").Value = -25472089.5

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 166669490
$ws.Range("I106").Value = 1000000000
$ws.Range("J106").Value = 3380
$ws.Range("K106").Value = 1000000000
$ws.Range("L106").Value = 3380
$ws.Range("M106").Value = -999999369
$ws.Range("N106").Value = -4642

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1315.1904
$ws.Range("I137").Value = 1229.875
$ws.Range("J137").Value = 1428.9445
$ws.Range("K137").Value = 3689.625
$ws.Range("L137").Value = 4286.833500000001
$ws.Range("M137").Value = -1139.625
$ws.Range("N137").Value = -9386.833500000001

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1524.5238
$ws.Range("I138").Value = 829.7778
$ws.Range("J138").Value = 3261.389
$ws.Range("K138").Value = 2489.3334
$ws.Range("L138").Value = 9784.167000000001
$ws.Range("M138").Value = 2650.6666
$ws.Range("N138").Value = -20064.167

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1339.2903
$ws.Range("I32").Value = 1250.159
$ws.Range("K32").Value = 1250.159
$ws.Range("M32").Value = -963.1590000000001

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 17597532
$ws.Range("I45").Value = 25718870
$ws.Range("J45").Value = 1299.5
$ws.Range("K45").Value = 25718870
$ws.Range("L45").Value = 1299.5
$ws.Range("M45").Value = -25718493
$ws.Range("N45").Value = -2053.5

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1020.9787
$ws.Range("I97").Value = 733.1667
$ws.Range("J97").Value = 1528.8823
$ws.Range("K97").Value = 733.1667
$ws.Range("L97").Value = 1528.8823
$ws.Range("M97").Value = -237.1667
$ws.Range("N97").Value = -2520.8823

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1045.5
$ws.Range("I122").Value = 798.55554
$ws.Range("J122").Value = 1786.3334
$ws.Range("K122").Value = 2395.66662
$ws.Range("L122").Value = 5359.0002
$ws.Range("M122").Value = 54.33338000000003
$ws.Range("N122").Value = -10259.0002

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 586.6087
$ws.Range("I94").Value = 504.8421
$ws.Range("K94").Value = 504.8421
$ws.Range("M94").Value = -53.84210000000002

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 71429600
$ws.Range("I107").Value = 125000580
$ws.Range("J107").Value = 1637.6666
$ws.Range("K107").Value = 125000580
$ws.Range("L107").Value = 1637.6666
$ws.Range("M107").Value = -124998660
$ws.Range("N107").Value = -5477.6666

# BSM row 109
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 19085
$ws.Range("J109").Value = 19085
$ws.Range("L109").Value = 19085
$ws.Range("N109").Value = -21859

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2587095.8
$ws.Range("I134").Value = 570.4666999999999
$ws.Range("J134").Value = 8556000
$ws.Range("K134").Value = 1711.4001
$ws.Range("L134").Value = 25668000
$ws.Range("M134").Value = 823.5999000000002
$ws.Range("N134").Value = -25673070

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 166668830
$ws.Range("I99").Value = 200001700
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 200001700
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -200000202
$ws.Range("N99").Value = -7496

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 166668830
$ws.Range("I126").Value = 200001700
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 600005100
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -600002630
$ws.Range("N126").Value = -18440

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 17545386
$ws.Range("I132").Value = 1597.4286
$ws.Range("J132").Value = 66667996
$ws.Range("K132").Value = 4792.2858
$ws.Range("L132").Value = 200003988
$ws.Range("M132").Value = -2262.2858
$ws.Range("N132").Value = -200009048

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1447.2821
$ws.Range("I134").Value = 1319.3914
$ws.Range("J134").Value = 1631.125
$ws.Range("K134").Value = 3958.1742
$ws.Range("L134").Value = 4893.375
$ws.Range("M134").Value = -1423.1742
$ws.Range("N134").Value = -9963.375

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 47622120
$ws.Range("I133").Value = 55558270
$ws.Range("J133").Value = 5200
$ws.Range("K133").Value = 166674810
$ws.Range("L133").Value = 15600
$ws.Range("M133").Value = -166669750
$ws.Range("N133").Value = -25720

# GSM row 24
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 2000000
$ws.Range("I24").Value = 2000000
$ws.Range("K24").Value = 2000000
$ws.Range("M24").Value = -1999827

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1870.6666
$ws.Range("I102").Value = 1306
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1306
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 316
$ws.Range("N102").Value = -6244

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2998.4167
$ws.Range("J126").Value = 3088.0908
$ws.Range("L126").Value = 9264.2724
$ws.Range("N126").Value = -14204.2724

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6989.5
$ws.Range("I132").Value = 4401.423
$ws.Range("J132").Value = 18204.5
$ws.Range("K132").Value = 13204.269
$ws.Range("L132").Value = 54613.5
$ws.Range("M132").Value = -10674.269
$ws.Range("N132").Value = -59673.5

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1546.3334
$ws.Range("I7").Value = 1434.6666
$ws.Range("K7").Value = 1434.6666
$ws.Range("M7").Value = -1322.6666

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 50002300
$ws.Range("I40").Value = 2750
$ws.Range("K40").Value = 2750
$ws.Range("M40").Value = -2614

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 33334170
$ws.Range("I93").Value = 47619756
$ws.Range("J93").Value = 1133.3334
$ws.Range("K93").Value = 47619756
$ws.Range("L93").Value = 1133.3334
$ws.Range("M93").Value = -47618508
$ws.Range("N93").Value = -3629.3334

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2708.3
$ws.Range("I100").Value = 3350
$ws.Range("J100").Value = 2609.577
$ws.Range("K100").Value = 3350
$ws.Range("L100").Value = 2609.577
$ws.Range("M100").Value = -2809
$ws.Range("N100").Value = -3691.577

# LTW row 117
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H117").Value = 23392
$ws.Range("J117").Value = 23392
$ws.Range("L117").Value = 23392
$ws.Range("N117").Value = -32570

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1546.3334
$ws.Range("I126").Value = 1434.6666
$ws.Range("K126").Value = 4303.9998
$ws.Range("M126").Value = -1833.9998

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 24317954
$ws.Range("I136").Value = 4466571.5
$ws.Range("J136").Value = 66667572
$ws.Range("K136").Value = 13399714.5
$ws.Range("L136").Value = 200002716
$ws.Range("M136").Value = -13397164.5
$ws.Range("N136").Value = -200007816

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 43479004
$ws.Range("I100").Value = 47619604
$ws.Range("J100").Value = 2666.5
$ws.Range("K100").Value = 95239208
$ws.Range("L100").Value = 5333
$ws.Range("M100").Value = -95238667
$ws.Range("N100").Value = -6415

# WVR row 118
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 21464
$ws.Range("J118").Value = 21464
$ws.Range("L118").Value = 21464
$ws.Range("N118").Value = -24778

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1501.7273
$ws.Range("I126").Value = 1052.3334
$ws.Range("J126").Value = 2041
$ws.Range("K126").Value = 3157.0002
$ws.Range("L126").Value = 6123
$ws.Range("M126").Value = -687.0001999999999
$ws.Range("N126").Value = -11063

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3597.5442
$ws.Range("I136").Value = 4568.289
$ws.Range("J136").Value = 1698.2609
$ws.Range("K136").Value = 13704.867
$ws.Range("L136").Value = 5094.7827
$ws.Range("M136").Value = -11154.867
$ws.Range("N136").Value = -10194.7827
